$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for a2cc017e-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-22 20:45:17"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for a2cc017e-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-22 20:45:07"
$wsZhCn.Range("K4").Value = "2016-08-22 20:45:32"

# de-de sheet: Latest HO Xliff Generate Date (Correspond Handoff Datetime) and
# Correspond Handback DateTime for a2cc017e-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-22 20:45:17"
$wsDeDe.Range("K4").Value = "2016-08-22 20:45:39"
